# Incorporate updated data from upstream processes through 2024
#
# Updates four data cells on Sheet1 (the chart's source range) to their
# new upstream values:
#   C26 (2024 / Energy Storage)        : 1019.68  -> 1034.68
#   E23 (2021 / Solar)                 : 975.159  -> 971.859
#   E25 (2023 / Solar)                 : 1872.644 -> 1877.544
#   E26 (2024 / Solar)                 : 2108.764 -> 2571.994

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26").Value = 1034.68
$ws.Range("E23").Value = 971.859
$ws.Range("E25").Value = 1877.544
$ws.Range("E26").Value = 2571.994

$wb.Save()
